# Insert a new data row at row 142 (pushing existing rows 142:233 down to 143:234)
# and populate it with the new "Brócoli" record for Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 142; this shifts rows 142-233 down to 143-234
# and copies formatting from the row below (matching the original author's edit).
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with values.
$ws.Range("A142").Value2 = 4
$ws.Range("B142").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C142").Value2 = "Los Lagos"
$ws.Range("D142").Value2 = 44518
$ws.Range("E142").Value2 = 10
$ws.Range("F142").Value2 = 100112023
$ws.Range("G142").Value2 = "Brócoli"
$ws.Range("H142").Value2 = "Sin especificar"
$ws.Range("I142").Value2 = "Primera"
$ws.Range("J142").Value2 = 500
$ws.Range("K142").Value2 = 1200
$ws.Range("L142").Value2 = 1200
$ws.Range("M142").Value2 = 1200
$ws.Range("N142").Value2 = "`$/unidad"
$ws.Range("O142").Value2 = "Región Metropolitana"
$ws.Range("P142").Value2 = 1200
$ws.Range("Q142").Value2 = 1
$ws.Range("R142").Value2 = "Hortaliza"

# Make sure the D142 cell keeps the date number format (style) used by the rest
# of column D, in case the insert did not carry it over.
$ws.Range("D142").NumberFormat = "YYYY-MM-DD HH:MM:SS"
